$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Emma Myers name to "Emma " (trailing space preserved)
$ws.Range("A4").Value = "Emma "

# Add new row 5 with example worker data
$ws.Range("A5").Value = "Trabajador de ejemplo"
$ws.Range("B5").Value = "descripcion`nde`nejemplo"
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 600
$ws.Range("E5").Value = "hola como estas"
$ws.Range("F5").Value = -123
$ws.Range("G5").Value = "mu bien"
$ws.Range("H5").Value = 456

# Multi-line text entry auto-expands row height; restore default auto height
$ws.Rows.Item(5).AutoFit()
